$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 13: assistantExaminer ----
$ws.Range("A13").Value = "assistantExaminer"

$ws.Range("B13").Value = "0947241"

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "kurnool_eGov@123"
$ws.Hyperlinks.Add($ws.Range("C13"), "mailto:kurnool_eGov@123", "", "", "kurnool_eGov@123")
# Hyperlinks.Add auto-applies the built-in "hyperlink" look (blue + underline).
# Put the font back the way the rest of column C looks so C13 matches C7:C12.
$ws.Range("C13").Font.Name = "Calibri"
$ws.Range("C13").Font.Size = 12
$ws.Range("C13").Font.ColorIndex = 1
$ws.Range("C13").Font.Underline = -4142
$ws.Range("C13").NumberFormat = "@"

$ws.Range("D13").Formula = "=FALSE()"
$ws.Range("D13").NumberFormat = '"TRUE";"TRUE";"FALSE"'

$ws.Rows.Item(13).RowHeight = 15.7

# ---- Row 14: examiner ----
$ws.Range("A14").Value = "examiner"

$ws.Range("B14").Value = "0800129"

$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "kurnool_eGov@123"
$ws.Hyperlinks.Add($ws.Range("C14"), "mailto:kurnool_eGov@123", "", "", "kurnool_eGov@123")
$ws.Range("C14").Font.Name = "Calibri"
$ws.Range("C14").Font.Size = 12
$ws.Range("C14").Font.ColorIndex = 1
$ws.Range("C14").Font.Underline = -4142
$ws.Range("C14").NumberFormat = "@"

$ws.Range("D14").Formula = "=FALSE()"

$ws.Rows.Item(14).RowHeight = 15.7

# ---- Selection, matching the post-edit cursor position ----
$ws.Range("D15").Select()
